# Applies the proofing-pass restructuring described by the commit diff:
# Word re-split runs and inserted <w:proofErr> spellStart/spellEnd markers
# around out-of-dictionary (English) words across the vocab list and the
# RBG / Mike Pence notes, plus a new empty "en-GB" paragraph.

$d = $word.ActiveDocument

function Replace-ParaRangeXml([int]$firstPara, [int]$lastPara, [string]$bodyXml) {
    $startPos = $d.Paragraphs($firstPara).Range.Start
    $endPos = $d.Paragraphs($lastPara).Range.End
    $rng = $d.Range($startPos, $endPos)
    $pkg = $bodyXml
    $rng.InsertXML($pkg)
}

# --- Vocabulary list (paragraphs 1-12) + new blank "en-GB" paragraph ---
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">To </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>be</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>deceiving</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>candlelight</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> vigil = tribut </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>with</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>candles</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>generaly</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> silent </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gathering</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) = un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>homage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> avec des bougies</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>saying</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>proverb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Drowsy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">= </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sleepy</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Defacement</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vandalism</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Shameless</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = sans honte</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Skull</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = crâne</w:t></w:r></w:p><w:p><w:r><w:t>Bold = brave/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>courageous</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bold</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en gras</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Bald = no </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hair</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>your</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>head</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wake-up</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> call = a warning </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cry</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">To </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>be</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shambles</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ruins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = être en ruine</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">A Coroner = a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>forensic</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>doctor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>legal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>doctor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 1 13 $xml1

# --- Paragraph 16: 2nd woman on the S.C. ---
$xml16 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="2410"/></w:pPr><w:r><w:t>2</w:t></w:r><w:r w:rsidRPr="00FA350E"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>woman</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> on the S.C.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 16 16 $xml16

# --- Paragraph 17: Fought for : (list item) ---
$xml17 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="2410"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Fought</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for : </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 17 17 $xml17

# --- Paragraph 18: Equality (list item) ---
$xml18 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="4"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Equality</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 18 18 $xml18

# --- Paragraph 19: Women's rights (abortion) (list item) ---
$xml19 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="4"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Women’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rights</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (abortion)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 19 19 $xml19

# --- Paragraph 20: Created proverb ... Ruth ! ---
$xml20 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Created</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>proverb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> : You </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>can’t</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>spell</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>truth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>without</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Ruth !</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 20 20 $xml20

# --- Paragraph 22: Segregation -> Jim Crow Laws ---
$xml22 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Segregation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> Jim Crow Laws</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 22 22 $xml22

# --- Paragraph 23: S.C. decision 1897 -> Seperate but equal ---
$xml23 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">S.C. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>decision</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 1897 </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> « </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Seperate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> but </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>equal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t> »</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 23 23 $xml23

# --- Paragraph 24: 1954 : Separation is inherently unequal ---
$xml24 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>1954 : « </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Separation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>is</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>inherently</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>unequal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t> »</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 24 24 $xml24

# --- Paragraph 26: Mike Pence is an extremist : he finds ... ---
$xml26 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Mike Pence </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>is</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> an </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>extremist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> : </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>he</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>finds</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>it</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> immoral to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>be</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>same</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> room </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>with</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>woman</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParaRangeXml 26 26 $xml26

"done"